# "changes in address and add"
#
# Resolves the TODO note in E21 ("need to change to city") by renaming the
# B21 column entry from "District" to "City", and fills in several missing
# API-method / endpoint-name cells (H6, F7, F19, F34, G34) in the "User"
# and "Product" sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) E21 held a reviewer note ("need to change to city"); the note has been
#    actioned, so clear it (format reverts to the plain unbordered-fill
#    look used by the rest of column E, copied from E20).
# ---------------------------------------------------------------------
$ws.Range("E20").Copy() | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").ClearContents()

# ---------------------------------------------------------------------
# 2) B21: "District" -> "City"
# ---------------------------------------------------------------------
$ws.Range("B21").Value = "City"

# ---------------------------------------------------------------------
# 3) Fill in the newly documented API entry points / methods. These cells
#    were blank; give them the same red-fill "placeholder" style already
#    used for the empty divider rows (e.g. G16), then set their text.
# ---------------------------------------------------------------------
$ws.Range("G16").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$ws.Range("F19").Value = "Update Adress"

$ws.Range("G16").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null
$ws.Range("F7").Value = "Update User"

$ws.Range("G16").Copy() | Out-Null
$ws.Range("F34").PasteSpecial(-4122) | Out-Null
$ws.Range("F34").Value = "Add A product"

$ws.Range("G16").Copy() | Out-Null
$ws.Range("G34").PasteSpecial(-4122) | Out-Null
$ws.Range("G34").Value = " need to implement in frontend"

# ---------------------------------------------------------------------
# 4) H6 was blank; it should carry the HTTP method ("GET") like the other
#    cells in column H. Match the look of the cell above it (H5).
# ---------------------------------------------------------------------
$ws.Range("H5").Copy() | Out-Null
$ws.Range("H6").PasteSpecial(-4122) | Out-Null
$ws.Range("H6").Value = "GET"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5) Leave the selection on H36, matching where editing ended up.
# ---------------------------------------------------------------------
$ws.Range("H36").Select() | Out-Null
